$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new logging entries for motor control gains (mailbox row 15 -> 12)
$ws.Range("C15").Value = "Output"
$ws.Range("D15").Value = "HO_VectorBox/motorControl"

# Update the active selection cell as recorded at save time
$ws.Range("J16").Select()
